$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Organograma")

# "Superintendente de T.I." -> "Superintendente de TI" (row 4 - Comitê Executivo)
$ws.Range("D4").Value = "Superintendente de TI"

# "TI" (Área) -> "Tecnologia da Informação" across the whole sheet
$ws.Range("E4").Value = "Tecnologia da Informação"
$ws.Range("E8").Value = "Tecnologia da Informação"
$ws.Range("E9").Value = "Tecnologia da Informação"
$ws.Range("E13").Value = "Tecnologia da Informação"
$ws.Range("E14").Value = "Tecnologia da Informação"
$ws.Range("E15").Value = "Tecnologia da Informação"
$ws.Range("E16").Value = "Tecnologia da Informação"
$ws.Range("E18").Value = "Tecnologia da Informação"
$ws.Range("E19").Value = "Tecnologia da Informação"
$ws.Range("E20").Value = "Tecnologia da Informação"

# "Analista" -> "Analista de Sistemas " (rows 18-20, Time Projeto)
$ws.Range("B18").Value = "Analista de Sistemas "
$ws.Range("B19").Value = "Analista de Sistemas "
$ws.Range("B20").Value = "Analista de Sistemas "

# Update the frozen pane / active selection to match the latest view state
$ws.Application.ActiveWindow.SplitRow = 2
$ws.Range("A15").Select()
$ws.Range("C25").Select()
